# Add a new "Profile47" test case row (row 48) to the Test Cases sheet,
# mirroring the formatting of the last existing row (row 47).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 47's formatting (styles/borders/fills) down into the new row 48,
# since Copy(Destination) duplicates both values and formatting.
$ws.Range("A47:E47").Copy($ws.Range("A48:E48"))

# Column C in the new row should have no explicit formatting (matches the
# target workbook, where C48 carries no style override).
$ws.Range("C48").ClearFormats()

# Fill in the new test case's data. Order matters so that new shared
# strings are appended in the same sequence as the target workbook.
$ws.Range("A48").Value = "Profile47"
$ws.Range("C48").Value = "Verify that profile call to Action(CTA) in a white box is getting  displayed in Summary field when your profile summary is blank"
$ws.Range("B48").Value = "OPQA-3323"
$ws.Range("D48").Value = "Y"

# Update the current selection to match the saved workbook state.
$null = $ws.Range("C39").Select()
